# checklist_forStudent4.1.xlsx — "setEadmin and seek help"
#
# Marks the "Sets up E-Admin and senior E-Admin" (D3) and "Seek help" (D28)
# checklist rows as done, using the same text/format already used elsewhere
# in the sheet for a completed row (copy the format from D5, which already
# carries that style), then writes the new status text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D5 already uses the "done" number format/font for this column; clone its
# formatting onto D3 and D28 before overwriting their values.
$ws.Range("D5").Copy()
$ws.Range("D3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D5").Copy()
$ws.Range("D28").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("D3").Value = "3done"
$ws.Range("D28").Value = "3done"

# Leave the cursor where the author left it when they saved the file.
$ws.Range("D4").Select() | Out-Null
